$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$oldVersion = "mines - January 30 (built on $oldStamp)"
$newVersion = "mines - January 30 (built on $newStamp)"

# --- "About" sheet ---
$about = $wb.Worksheets.Item("About")

$about.Range("A2").Value = "Version: $newVersion"
$about.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Huashan Coal Mine, China, M1965, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$data = $wb.Worksheets.Item("Boundaries and methane sources")

$lastRow = $data.Cells.Item($data.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $data.Cells.Item($r, 19)  # column S
    if ($cell.Value2 -ne $null -and $cell.Value2 -like "mines - January 30 (built on *") {
        $cell.Value = $newVersion
    }
}
